$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $isNum = $Text -match '^[+-]?[0-9]+(\.[0-9]+)?$'
    if ($isNum) {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

Set-TextValue $ws.Range('D2') '63.533.70'
Set-TextValue $ws.Range('E2') '  +5.84%  '
Set-TextValue $ws.Range('D3') '3.399.09'
Set-TextValue $ws.Range('E3') '  +6.66%  '
Set-TextValue $ws.Range('E4') '  +0.01%  '
Set-TextValue $ws.Range('D5') '578.31'
Set-TextValue $ws.Range('E5') '  +7.82%  '
Set-TextValue $ws.Range('D6') '155.89'
Set-TextValue $ws.Range('E6') '  +7.48%  '
Set-TextValue $ws.Range('D7') '1.00'
Set-TextValue $ws.Range('E7') '  -0.02%  '
Set-TextValue $ws.Range('D8') '3.407.27'
Set-TextValue $ws.Range('E8') '  +6.75%  '
Set-TextValue $ws.Range('D9') '0.532'
Set-TextValue $ws.Range('E10') '  +2.45%  '
Set-TextValue $ws.Range('E11') '  +7.49%  '
Set-TextValue $ws.Range('D12') '0.437'
Set-TextValue $ws.Range('E12') '  +1.24%  '
Set-TextValue $ws.Range('D13') '3.982.32'
Set-TextValue $ws.Range('E13') '  +6.55%  '
Set-TextValue $ws.Range('E14') '  +0.42%  '
Set-TextValue $ws.Range('E15') '  +7.09%  '
Set-TextValue $ws.Range('D16') '27.08'
Set-TextValue $ws.Range('E16') '  +4.98%  '
Set-TextValue $ws.Range('D17') '63.594.20'
Set-TextValue $ws.Range('E17') '  +5.92%  '
Set-TextValue $ws.Range('D18') '3.388.75'
Set-TextValue $ws.Range('E18') '  +6.24%  '
Set-TextValue $ws.Range('D19') '6.37'
Set-TextValue $ws.Range('E19') '  +2.19%  '
Set-TextValue $ws.Range('D20') '14.05'
Set-TextValue $ws.Range('E20') '  +6.24%  '
Set-TextValue $ws.Range('D21') '8.47'
Set-TextValue $ws.Range('E21') '  +3.47%  '
Set-TextValue $ws.Range('D22') '387.06'
Set-TextValue $ws.Range('E22') '  +4.95%  '
Set-TextValue $ws.Range('D23') '0.999'
Set-TextValue $ws.Range('D24') '0.535'
Set-TextValue $ws.Range('E24') '  +2.46%  '
Set-TextValue $ws.Range('D25') '70.94'
Set-TextValue $ws.Range('E25') '  +2.14%  '
Set-TextValue $ws.Range('D26') '9.56'
Set-TextValue $ws.Range('E26') '  +11.43%  '
Set-TextValue $ws.Range('E27') '  +6.90%  '
Set-TextValue $ws.Range('E28') '  +17.93%  '
Set-TextValue $ws.Range('D29') '1.00'
Set-TextValue $ws.Range('E29') '  +1.10%  '
Set-TextValue $ws.Range('E30') '  +7.86%  '
Set-TextValue $ws.Range('D31') '6.60'
Set-TextValue $ws.Range('E31') '  +8.00%  '
Set-TextValue $ws.Range('D32') '1.36'
Set-TextValue $ws.Range('E32') '  +14.41%  '
Set-TextValue $ws.Range('D33') '5.64'
Set-TextValue $ws.Range('E33') '  +7.14%  '
Set-TextValue $ws.Range('D34') '23.13'
Set-TextValue $ws.Range('E34') '  +2.95%  '
Set-TextValue $ws.Range('B35') 'USDe'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D35') '0.997'
Set-TextValue $ws.Range('E35') '  -0.10%  '
Set-TextValue $ws.Range('B36') 'Aptos'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D36') '6.69'
Set-TextValue $ws.Range('E36') '  +1.75%  '
Set-TextValue $ws.Range('B37') 'ImmutableX'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D37') '1.49'
Set-TextValue $ws.Range('E37') '  +9.82%  '
Set-TextValue $ws.Range('B38') 'Monero'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D38') '158.09'
Set-TextValue $ws.Range('E38') '  +0.30%  '
Set-TextValue $ws.Range('B39') 'Stacks'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D39') '1.87'
Set-TextValue $ws.Range('E39') '  +10.95%  '
Set-TextValue $ws.Range('B40') 'EnergySwap'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D40') '27.52'
Set-TextValue $ws.Range('E40') '  +4.17%  '
Set-TextValue $ws.Range('B41') 'Hedera'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D41') '0.0763'
Set-TextValue $ws.Range('E41') '  +7.86%  '
Set-TextValue $ws.Range('B42') 'Maker'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D42') '2.906.15'
Set-TextValue $ws.Range('E42') '  +4.18%  '
Set-TextValue $ws.Range('B43') 'VeChain'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D43') '0.0322'
Set-TextValue $ws.Range('E43') '  +4.76%  '
Set-TextValue $ws.Range('B44') 'Mantle'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D44') '0.764'
Set-TextValue $ws.Range('E44') '  +6.31%  '
Set-TextValue $ws.Range('B45') 'OKB'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D45') '41.46'
Set-TextValue $ws.Range('E45') '  +4.11%  '
Set-TextValue $ws.Range('B46') 'Filecoin'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D46') '4.31'
Set-TextValue $ws.Range('E46') '  +2.14%  '
Set-TextValue $ws.Range('B47') 'ONDO'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range('D47') '1.07'
Set-TextValue $ws.Range('E47') '  +8.51%  '
Set-TextValue $ws.Range('D48') '22.49'
Set-TextValue $ws.Range('E48') '  +9.33%  '
Set-TextValue $ws.Range('B49') 'RenzoRestakedETH'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue $ws.Range('D49') '3.444.08'
Set-TextValue $ws.Range('E49') '  +6.69%  '
Set-TextValue $ws.Range('B50') 'Bittensor'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D50') '299.23'
Set-TextValue $ws.Range('E50') '  +14.17%  '
Set-TextValue $ws.Range('B51') 'Stellar'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D51') '0.103'
Set-TextValue $ws.Range('E51') '  -2.03%  '
